# Apply small adjustments to DD for EPICP to reflect dataset file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Fix casing of the PAL variable name (row 7, column B)
$ws.Range("B7").Value = "pal"

# Append new variable rows (index 66-77 -> sheet rows 67-78)
$newRows = @(
    @(66, "VEGETABLES_02", "Vegetable intake [g/d]", "decimal"),
    @(67, "LEGUMES_TOT_03", "Total legumes intake [g/d]", "decimal"),
    @(68, "FRUITS_TOT_04", "Total fruit intake [g/d]", "decimal"),
    @(69, "RED_MEAT_0701", "Intake of red meat (mammals meat) [g/d]", "decimal"),
    @(70, "PROCMEAT_0704", "Intake of processed or preserved meat [g/d]", "decimal"),
    @(71, "SUGAR_CONFECT_11", "Intake of sugar and similar, confectionery and water-based sweet desserts [g/d]", "decimal"),
    @(72, "CAKES_12", "Intake of cakes and fine bakery products [g/d]", "decimal"),
    @(73, "FRUITVEG_JUICE_1301", "Intake of fruit and vegetable juices [g/d]", "decimal"),
    @(74, "SOFTDRINKS_1302", "Intake of soft drinks [g/d]", "decimal"),
    @(75, "COFFEE_130301", "Coffee intake [g/d]", "decimal"),
    @(76, "TEA_130302", "Tea intake [g/d]", "decimal"),
    @(77, "ART_SWEETENER_170201", "Intake of artificial sweeteners (e.g., aspartam, saccharine) [g/d]", "decimal")
)

$startRow = 67
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
}
